$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(717).Delete()
